$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new daily-scrape row (row 65) at the bottom of the sheet.
# Column A holds a date-like string that must stay literal text (matching
# the existing "2025/10/05" cells above it), so force the text number
# format before assigning the value to stop Excel auto-converting it to a
# date serial number.
$ws.Range("A65").NumberFormat = "@"
$ws.Range("A65").Value = "2025/10/05"
$ws.Range("B65").Value = "日"
$ws.Range("C65").Value = 20
$ws.Range("D65").Value = 56
